# GitHub Actions "Updated cryptos list" refresh.
# Re-writes the Price (D) / Volume(1h) (E) columns with the latest scraped
# values, and (for this run) OKB and TheGraph swapped ranking positions so
# their row 39 / row 40 data (Coin, Link, Price, Volume) are exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing it to stay Text-typed.
# The source data has numeric-looking strings (e.g. "134.10", "0.940")
# whose exact formatting (trailing zeros, thousands-separator dots, etc.)
# must be preserved, matching the original inlineStr cells. A plain
# Value assignment lets Excel silently reinterpret such strings as
# numbers (dropping formatting), so for values that look numeric we
# briefly mark the cell as Text ("@"), write the value, then clear the
# explicit formatting again so no stray number-format/style is left on
# the cell (it returns to the workbook default, unstyled state).
function Set-TextValue($range, $value) {
    if ($value -match '^\s*[-+]?\d+(\.\d+)?\s*$') {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") '66.602.00'
Set-TextValue $ws.Range("E2") '  -0.96%  '
Set-TextValue $ws.Range("D3") '3.087.01'
Set-TextValue $ws.Range("E3") '  -1.42%  '
Set-TextValue $ws.Range("E4") '  +0.04%  '
Set-TextValue $ws.Range("D5") '575.91'
Set-TextValue $ws.Range("E5") '  -0.94%  '
Set-TextValue $ws.Range("D6") '172.20'
Set-TextValue $ws.Range("E6") '  -1.31%  '
Set-TextValue $ws.Range("E7") '  +0.07%  '
Set-TextValue $ws.Range("D8") '3.085.48'
Set-TextValue $ws.Range("E8") '  -1.39%  '
Set-TextValue $ws.Range("E9") '  -1.81%  '
Set-TextValue $ws.Range("D10") '6.36'
Set-TextValue $ws.Range("E10") '  -0.93%  '
Set-TextValue $ws.Range("E11") '  -2.75%  '
Set-TextValue $ws.Range("E12") '  -2.50%  '
Set-TextValue $ws.Range("E13") '  -4.40%  '
Set-TextValue $ws.Range("D14") '35.64'
Set-TextValue $ws.Range("E14") '  -5.15%  '
Set-TextValue $ws.Range("E15") '  -0.84%  '
Set-TextValue $ws.Range("D16") '3.601.67'
Set-TextValue $ws.Range("E16") '  -1.31%  '
Set-TextValue $ws.Range("D17") '66.549.03'
Set-TextValue $ws.Range("E17") '  -0.90%  '
Set-TextValue $ws.Range("E18") '  -2.96%  '
Set-TextValue $ws.Range("D19") '16.79'
Set-TextValue $ws.Range("E19") '  +2.40%  '
Set-TextValue $ws.Range("D20") '3.088.77'
Set-TextValue $ws.Range("E20") '  -1.26%  '
Set-TextValue $ws.Range("D21") '483.93'
Set-TextValue $ws.Range("E21") '  -2.01%  '
Set-TextValue $ws.Range("E22") '  -1.94%  '
Set-TextValue $ws.Range("E23") '  -3.46%  '
Set-TextValue $ws.Range("E24") '  -1.32%  '
Set-TextValue $ws.Range("D25") '12.62'
Set-TextValue $ws.Range("E25") '  -5.62%  '
Set-TextValue $ws.Range("E26") '  -3.57%  '
Set-TextValue $ws.Range("D27") '10.04'
Set-TextValue $ws.Range("E27") '  -4.03%  '
Set-TextValue $ws.Range("E28") '  -0.07%  '
Set-TextValue $ws.Range("D29") '7.96'
Set-TextValue $ws.Range("E29") '  +0.05%  '
Set-TextValue $ws.Range("E30") '  -4.77%  '
Set-TextValue $ws.Range("E31") '  -4.22%  '
Set-TextValue $ws.Range("E32") '  -3.06%  '
Set-TextValue $ws.Range("E33") '  -2.97%  '
Set-TextValue $ws.Range("D34") '0.0₃0934'
Set-TextValue $ws.Range("E34") '  -1.66%  '
Set-TextValue $ws.Range("E35") '  +0.13%  '
Set-TextValue $ws.Range("D36") '48.33'
Set-TextValue $ws.Range("E36") '  +3.24%  '
Set-TextValue $ws.Range("D37") '5.57'
Set-TextValue $ws.Range("E37") '  -5.98%  '
Set-TextValue $ws.Range("D38") '0.940'
Set-TextValue $ws.Range("E38") '  -3.94%  '
Set-TextValue $ws.Range("B39") 'OKB'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D39") '48.98'
Set-TextValue $ws.Range("E39") '  -2.40%  '
Set-TextValue $ws.Range("B40") 'TheGraph'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D40") '0.309'
Set-TextValue $ws.Range("E40") '  -1.44%  '
Set-TextValue $ws.Range("E42") '  -5.40%  '
Set-TextValue $ws.Range("D44") '2.59'
Set-TextValue $ws.Range("E44") '  -0.93%  '
Set-TextValue $ws.Range("D45") '2.774.51'
Set-TextValue $ws.Range("E45") '  -2.37%  '
Set-TextValue $ws.Range("E46") '  -2.78%  '
Set-TextValue $ws.Range("D47") '367.61'
Set-TextValue $ws.Range("E47") '  -5.06%  '
Set-TextValue $ws.Range("D48") '134.10'
Set-TextValue $ws.Range("D50") '24.38'
Set-TextValue $ws.Range("E50") '  -2.76%  '
Set-TextValue $ws.Range("D51") '2.16'
Set-TextValue $ws.Range("E51") '  -3.02%  '
